# Auto-generated edit script applying Behemoth_Profits.xlsx diff
# Updates currentAveragePrice / LevePrice / LeveProfit figures across the ALC,
# ARM, BSM, CRP, CUL, GSM, LTW, and WVR sheets to reflect refreshed market data.

$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 25050
$ws.Range("I20").Value = 25050
$ws.Range("K20").Value = 25050
$ws.Range("M20").Value = -24820
$ws.Range("H35").Value = 25050
$ws.Range("I35").Value = 25050
$ws.Range("K35").Value = 25050
$ws.Range("M35").Value = -24671
$ws.Range("H40").Value = 3871.842
$ws.Range("I40").Value = 2403
$ws.Range("J40").Value = 4396.4287
$ws.Range("K40").Value = 2403
$ws.Range("L40").Value = 4396.4287
$ws.Range("M40").Value = -2228
$ws.Range("N40").Value = -4746.4287
$ws.Range("H86").Value = 9545.046
$ws.Range("I86").Value = 9599.6
$ws.Range("K86").Value = 9599.6
$ws.Range("M86").Value = -8476.6
$ws.Range("H89").Value = 9545.046
$ws.Range("I89").Value = 9599.6
$ws.Range("K89").Value = 47998
$ws.Range("M89").Value = -42382
$ws.Range("H99").Value = 1707.5333
$ws.Range("I99").Value = 443.25
$ws.Range("J99").Value = 6764.6665
$ws.Range("K99").Value = 1329.75
$ws.Range("L99").Value = 20293.9995
$ws.Range("M99").Value = 168.25
$ws.Range("N99").Value = -23289.9995
$ws.Range("H112").Value = 4800
$ws.Range("J112").Value = 5750
$ws.Range("L112").Value = 17250
$ws.Range("N112").Value = -19466
$ws.Range("H118").Value = 1545.4445
$ws.Range("I118").Value = 407.5
$ws.Range("J118").Value = 2455.8
$ws.Range("K118").Value = 1222.5
$ws.Range("L118").Value = 7367.400000000001
$ws.Range("M118").Value = 434.5
$ws.Range("N118").Value = -10681.4
$ws.Range("H137").Value = 8357.333000000001
$ws.Range("I137").Value = 2969.2
$ws.Range("K137").Value = 8907.599999999999
$ws.Range("M137").Value = -6357.599999999999
$ws.Range("H138").Value = 1838.4193
$ws.Range("I138").Value = 1057.579
$ws.Range("J138").Value = 3074.75
$ws.Range("K138").Value = 3172.737
$ws.Range("L138").Value = 9224.25
$ws.Range("M138").Value = 1967.263
$ws.Range("N138").Value = -19504.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9092135
$ws.Range("I32").Value = 9092135
$ws.Range("K32").Value = 9092135
$ws.Range("M32").Value = -9091848
$ws.Range("H74").Value = 13168565
$ws.Range("I74").Value = 20834022
$ws.Range("K74").Value = 20834022
$ws.Range("M74").Value = -20833148
$ws.Range("H77").Value = 13168565
$ws.Range("I77").Value = 20834022
$ws.Range("K77").Value = 104170110
$ws.Range("M77").Value = -104165742
$ws.Range("H88").Value = 2890.1
$ws.Range("I88").Value = 2498.8
$ws.Range("J88").Value = 3281.4
$ws.Range("K88").Value = 2498.8
$ws.Range("L88").Value = 3281.4
$ws.Range("M88").Value = -2092.8
$ws.Range("N88").Value = -4093.4
$ws.Range("H91").Value = 2890.1
$ws.Range("I91").Value = 2498.8
$ws.Range("J91").Value = 3281.4
$ws.Range("K91").Value = 2498.8
$ws.Range("L91").Value = 3281.4
$ws.Range("M91").Value = -1094.8
$ws.Range("N91").Value = -6089.4
$ws.Range("H97").Value = 898.6957
$ws.Range("I97").Value = 910
$ws.Range("K97").Value = 910
$ws.Range("M97").Value = -414
$ws.Range("H102").Value = 13256.963
$ws.Range("I102").Value = 18242.23
$ws.Range("K102").Value = 18242.23
$ws.Range("M102").Value = -16620.23
$ws.Range("H132").Value = 6068.879
$ws.Range("I132").Value = 4426.0713
$ws.Range("K132").Value = 13278.2139
$ws.Range("M132").Value = -10748.2139

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4888.357
$ws.Range("I20").Value = 4649.1665
$ws.Range("J20").Value = 6323.5
$ws.Range("K20").Value = 4649.1665
$ws.Range("L20").Value = 6323.5
$ws.Range("M20").Value = -4402.1665
$ws.Range("N20").Value = -6817.5
$ws.Range("H63").Value = 116999.5
$ws.Range("J63").Value = 116999.5
$ws.Range("L63").Value = 116999.5
$ws.Range("N63").Value = -118371.5
$ws.Range("H66").Value = 116999.5
$ws.Range("J66").Value = 116999.5
$ws.Range("L66").Value = 350998.5
$ws.Range("N66").Value = -357862.5
$ws.Range("H86").Value = 3708.7693
$ws.Range("I86").Value = 3246.7273
$ws.Range("J86").Value = 6250
$ws.Range("K86").Value = 3246.7273
$ws.Range("L86").Value = 6250
$ws.Range("M86").Value = -2123.7273
$ws.Range("N86").Value = -8496
$ws.Range("H89").Value = 3708.7693
$ws.Range("I89").Value = 3246.7273
$ws.Range("J89").Value = 6250
$ws.Range("K89").Value = 16233.6365
$ws.Range("L89").Value = 31250
$ws.Range("M89").Value = -10617.6365
$ws.Range("N89").Value = -42482
$ws.Range("H94").Value = 1308.871
$ws.Range("I94").Value = 1302.1666
$ws.Range("K94").Value = 1302.1666
$ws.Range("M94").Value = -851.1666
$ws.Range("H107").Value = 1321.5333
$ws.Range("I107").Value = 1123.8462
$ws.Range("J107").Value = 2606.5
$ws.Range("K107").Value = 1123.8462
$ws.Range("L107").Value = 2606.5
$ws.Range("M107").Value = 796.1538
$ws.Range("N107").Value = -6446.5
$ws.Range("H134").Value = 29324.54
$ws.Range("I134").Value = 1409.5
$ws.Range("K134").Value = 4228.5
$ws.Range("M134").Value = -1693.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 943678.9399999999
$ws.Range("I31").Value = 27285
$ws.Range("J31").Value = 1233066.5
$ws.Range("K31").Value = 27285
$ws.Range("L31").Value = 1233066.5
$ws.Range("M31").Value = -26990
$ws.Range("N31").Value = -1233656.5
$ws.Range("H34").Value = 943678.9399999999
$ws.Range("I34").Value = 27285
$ws.Range("J34").Value = 1233066.5
$ws.Range("K34").Value = 27285
$ws.Range("L34").Value = 1233066.5
$ws.Range("M34").Value = -27083
$ws.Range("N34").Value = -1233470.5
$ws.Range("H134").Value = 531876.9399999999
$ws.Range("J134").Value = 11389.25
$ws.Range("L134").Value = 34167.75
$ws.Range("N134").Value = -39237.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 700
$ws.Range("I14").Value = 700
$ws.Range("K14").Value = 2100
$ws.Range("M14").Value = -1927
$ws.Range("H25").Value = 3334233.2
$ws.Range("J25").Value = 5001250
$ws.Range("L25").Value = 15003750
$ws.Range("N25").Value = -15004088
$ws.Range("H30").Value = 3334233.2
$ws.Range("J30").Value = 5001250
$ws.Range("L30").Value = 15003750
$ws.Range("N30").Value = -15003954
$ws.Range("H38").Value = 52.75
$ws.Range("I38").Value = 11
$ws.Range("J38").Value = 66.666664
$ws.Range("K38").Value = 33
$ws.Range("L38").Value = 199.999992
$ws.Range("M38").Value = 314
$ws.Range("N38").Value = -893.999992

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").Value = ""
$ws.Range("H70").Value = 4700
$ws.Range("I70").Value = 4700
$ws.Range("K70").Value = 4700
$ws.Range("M70").Value = -4430
$ws.Range("H73").Value = 4700
$ws.Range("I73").Value = 4700
$ws.Range("K73").Value = 4700
$ws.Range("M73").Value = -3764

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3906
$ws.Range("I46").Value = 3638.8
$ws.Range("K46").Value = 3638.8
$ws.Range("M46").Value = -3450.8
$ws.Range("H82").Value = 842.58826
$ws.Range("J82").Value = 1330.4
$ws.Range("L82").Value = 1330.4
$ws.Range("N82").Value = -2052.4
$ws.Range("H85").Value = 842.58826
$ws.Range("J85").Value = 1330.4
$ws.Range("L85").Value = 1330.4
$ws.Range("N85").Value = -3826.4
$ws.Range("H132").Value = 386756.12
$ws.Range("I132").Value = 477972.47
$ws.Range("J132").Value = 147313.25
$ws.Range("K132").Value = 1433917.41
$ws.Range("L132").Value = 441939.75
$ws.Range("M132").Value = -1431387.41
$ws.Range("N132").Value = -446999.75
$ws.Range("H136").Value = 57770.75
$ws.Range("I136").Value = 2869.5833
$ws.Range("K136").Value = 8608.749899999999
$ws.Range("M136").Value = -6058.749899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").Value = ""
$ws.Range("H136").Value = 10707.77
$ws.Range("I136").Value = 1389.8889
$ws.Range("K136").Value = 4169.6667
$ws.Range("M136").Value = -1619.6667
